$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: "Key " header, bold font
$ws.Range("A5").Value = "Key "

# Row 6-11: field name in column A, description in column B
$ws.Range("A6").Value = "date"
$ws.Range("B6").Value = "The date in 3-letter month + 2-number date (e.g. jun04, jan15, oct20, aug01... ). Important but inconvenient: excel tends to default to one of the data farmats. Counteract that by typing the apostrophe before text: 'jun04 "

$ws.Range("A7").Value = "box"
$ws.Range("B7").Value = "The numeric of the Firseting box."

$ws.Range("A8").Value = "channel"
$ws.Range("B8").Value = "The number of the channel (1-4)"

$ws.Range("A9").Value = "cycle"
$ws.Range("B9").Value = "Which measurement cycle will be cleaned? (cycle 2, cycle 3, cycle 4) Cycle is assumed to be MMR "

$ws.Range("A10").Value = "start"
$ws.Range("B10").Value = "indicate the start time (minutes) of the measurement section to be KEPT, write 0 to discard this cycle."

$ws.Range("A11").Value = "end"
$ws.Range("B11").Value = "indicate the end time of the section to be KEPT, write 0 to discard this cycle."

# Row 12 - empty styled cells (Arial font)
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""

# Fonts / Fills
$ws.Range("A5").Font.Bold = $true

$ws.Range("A6:B11").Interior.Color = 65535

$ws.Range("A12:B12").Font.Name = "Arial"

$ws.Range("A6").Select()
$ws.Range("A6:B11").Select()
